# ---------------------------------------------------------------------------
# "Added reports upload feature"
# Rebuilds the vocabulary sheet: replaces the placeholder rows (hhfhf/geeyy/hff)
# with real dictionary entries, adds a long-example row plus an inline
# "too-long field" note, pads the sheet with a few pre-formatted blank rows,
# and refreshes the sheet view (zoom, column widths, selection).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- clear the old placeholder rows (A2:B3) before rebuilding -------------
$ws.Range("A2:D3").ClearContents()

# ---- header row -------------------------------------------------------
$ws.Range("A1").Value = "Слово"
$ws.Range("B1").Value = "Мова"
$ws.Range("C1").Value = "Категорія"
$ws.Range("D1").Value = "Переклад"

# ---- data rows ----------------------------------------------------------
$ws.Range("A2").Value = "applicable"
$ws.Range("B2").Value = "Англійська"
$ws.Range("C2").Value = "Прикметник"
$ws.Range("D2").Value = "допустимий"

$ws.Range("A3").Value = "short"
$ws.Range("B3").Value = "Англійська"
$ws.Range("C3").Value = "Прикметник"
$ws.Range("D3").Value = "короткий"

$ws.Range("A4").Value = "Lorem ipsum dolor sit amet, consecterur adipiscing elit."
$ws.Range("B4").Value = "Латинська"
$ws.Range("C4").Value = "Фраза"
$ws.Range("D4").Value = "…"

$ws.Range("A5").Value = "unreachable"
$ws.Range("B5").Value = "Англійська"
$ws.Range("C5").Value = "Прикметник"
$ws.Range("D5").Value = "недосяжний"
$ws.Range("F5").Value = "№2."
$ws.Range("G5").Value = "Файл із задовгим полем у записі (>50 символів)"

$ws.Range("A6").Value = "record"
$ws.Range("B6").Value = "Англійська"
$ws.Range("C6").Value = "Іменник"
$ws.Range("D6").Value = "запис"

Write-Host "values done"
